# Update the "想去人数" (want-to-go count) figures in column F on both the
# "展览" and "全部类型" sheets, reflecting newer counts scraped at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 2344
    $ws.Range("F3").Value = 1861
    $ws.Range("F5").Value = 1137
    $ws.Range("F6").Value = 1117
    $ws.Range("F8").Value = 5956
}
